$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the Kyrgyz title text in A1 (shared string text correction).
#    This text replaces the old (slightly different) Kyrgyz translation; the
#    shared-strings table will automatically reshuffle so that the now-unused
#    old string slot gets reused by "-" and the corrected text is appended.
$ws.Range("A1").Value = "3.9.2 Коопсуздук суунун, коопсуздук санитариянын жана гигиенанын жоктугунан болгон өлүм"

# 2. Add a new "2022" column (S) with its data, copying number formatting
#    from the corresponding cell in the preceding "2021" column (R).
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 1.2

$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").Value = 2.7

$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Value = 0.9

$ws.Range("R8").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$ws.Range("S8").Value = 0.4

$ws.Range("R9").Copy()
$ws.Range("S9").PasteSpecial(-4122)
$ws.Range("S9").Value = 0.7

$ws.Range("R10").Copy()
$ws.Range("S10").PasteSpecial(-4122)
$ws.Range("S10").Value = 0.9

$ws.Range("R11").Copy()
$ws.Range("S11").PasteSpecial(-4122)
$ws.Range("S11").Value = 1.1000000000000001

$ws.Range("R12").Copy()
$ws.Range("S12").PasteSpecial(-4122)
$ws.Range("S12").Value = 2.7

$ws.Range("R13").Copy()
$ws.Range("S13").PasteSpecial(-4122)
$ws.Range("S13").Value = 0.4

$ws.Range("R14").Copy()
$ws.Range("S14").PasteSpecial(-4122)
$ws.Range("S14").Value = 0.6

# 3. Reset the selection away from the now out-of-range S17 cell back to A1.
$ws.Range("A1").Select()
